$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.265.94"
$ws.Range("E2").Value = "  +3.52%  "
$ws.Range("D3").Value = "1.920.00"
$ws.Range("E3").Value = "  +3.01%  "
$ws.Range("E4").Value = "  -1.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.09"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4855"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3852"
$ws.Range("E8").Value = "  +3.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07412"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9535"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.05"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07815"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "1.933.87"
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.561"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.668"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.37"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008909"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "28.263.69"
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.07"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.179"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").Value = "2.159.26"
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.51"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.59"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.133"
$ws.Range("E28").Value = "  +6.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.17"
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.053"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08920"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.362"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.254"
$ws.Range("E33").Value = "  +4.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7850"
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.693"
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.804"
$ws.Range("E36").Value = "  +4.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.134"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02058"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05391"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5598"
$ws.Range("E40").Value = "  +4.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.033"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.150"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.639"
$ws.Range("E43").Value = "  +2.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1540"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4950"
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.77"
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.02"
$ws.Range("E47").Value = "  +4.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.681"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.51"
$ws.Range("E50").Value = "  +4.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06151"
$ws.Range("E51").Value = "  +0.99%  "
